$wb2 = $excel.ActiveWorkbook
$ws  = $wb2.ActiveSheet

$nbsp = [char]0xA0

function Set-TextValue($cell, [string]$value) {
    # Force the cell to stay a TEXT value (Excel would otherwise silently
    # coerce pure-numeric-looking strings like "11" or "33.0" into numbers,
    # dropping trailing zeros / the text type). The leading apostrophe
    # forces text entry; resetting the cell Style back to "Normal"
    # afterwards clears the quotePrefix style Excel tacks on, so the cell
    # ends up on the same (default) style index as before the edit.
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# ---- Row 6 ----
Set-TextValue $ws.Cells.Item(6,2) "11"
$ws.Cells.Item(6,3).Value = "Randolph " + $nbsp
$ws.Cells.Item(6,4).Value = "Bridgette " + $nbsp
$ws.Cells.Item(6,5).Value = "-8.86,-5.32"
$ws.Cells.Item(6,6).Value = "Lenny(father): 0505536740"
Set-TextValue $ws.Cells.Item(6,8) "33.0"

# ---- Row 7 ----
Set-TextValue $ws.Cells.Item(7,2) "0"
$ws.Cells.Item(7,3).Value = "Trudie " + $nbsp
$ws.Cells.Item(7,4).Value = "Fleta " + $nbsp
$ws.Cells.Item(7,5).Value = "-5.2,-4.66"
$ws.Cells.Item(7,6).Value = "Anneliese(father): 0548973345"
Set-TextValue $ws.Cells.Item(7,8) "28.0"

# ---- Row 8 ----
Set-TextValue $ws.Cells.Item(8,2) "16"
$ws.Cells.Item(8,3).Value = "Collette " + $nbsp
$ws.Cells.Item(8,4).Value = "Billi " + $nbsp
$ws.Cells.Item(8,5).Value = "-4.8,-6.74"
$ws.Cells.Item(8,6).Value = "Elias(mother): 0578741979"
Set-TextValue $ws.Cells.Item(8,8) "24.0"

# ---- Row 9 ----
Set-TextValue $ws.Cells.Item(9,2) "7"
$ws.Cells.Item(9,3).Value = "Wyatt " + $nbsp
$ws.Cells.Item(9,4).Value = "Willette " + $nbsp
$ws.Cells.Item(9,5).Value = "-4.29,-7.75"
$ws.Cells.Item(9,6).Value = "Antionette(father): 0557331799"
$ws.Cells.Item(9,7).Value = "7:11:00"
Set-TextValue $ws.Cells.Item(9,8) "22.0"

# ---- Row 10 ----
Set-TextValue $ws.Cells.Item(10,2) "8"
$ws.Cells.Item(10,3).Value = "Marni " + $nbsp
$ws.Cells.Item(10,4).Value = "Shanika " + $nbsp
$ws.Cells.Item(10,5).Value = "-1.97,-7.93"
$ws.Cells.Item(10,6).Value = "Lady(mother): 0560804012"
$ws.Cells.Item(10,7).Value = "7:14:00"
Set-TextValue $ws.Cells.Item(10,8) "19.0"

# ---- Row 11 ----
Set-TextValue $ws.Cells.Item(11,2) "19"
$ws.Cells.Item(11,3).Value = "Jeanine " + $nbsp
$ws.Cells.Item(11,4).Value = "Janee " + $nbsp
$ws.Cells.Item(11,5).Value = "2.72,-7.47"
$ws.Cells.Item(11,6).Value = "Teresa(mother): 0517627420"
$ws.Cells.Item(11,7).Value = "7:21:00"
Set-TextValue $ws.Cells.Item(11,8) "12.0"

# ---- Row 12 ----
Set-TextValue $ws.Cells.Item(12,2) "3"
$ws.Cells.Item(12,3).Value = "Alexia " + $nbsp
$ws.Cells.Item(12,4).Value = "Ramonita " + $nbsp
$ws.Cells.Item(12,5).Value = "-0.99,0.61"
$ws.Cells.Item(12,6).Value = "Han(father): 0567537032"
$ws.Cells.Item(12,7).Value = "7:31:00"
Set-TextValue $ws.Cells.Item(12,8) "2.0"

# ---- Row 13 (school row) ----
$ws.Cells.Item(13,7).Value = "7:33:00"

# ---- Row 15 (time) ----
Set-TextValue $ws.Cells.Item(15,2) "33.0"
